$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "93.236.82"
$ws.Range("E2").Value = "  -3.03%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.313.37"
$ws.Range("E3").Value = "  -5.90%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "229.00"
$ws.Range("E5").Value = "  -6.62%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "616.54"
$ws.Range("E6").Value = "  -5.07%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.35"
$ws.Range("E7").Value = "  -4.88%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.382"
$ws.Range("E8").Value = "  -7.94%  "

$ws.Range("E9").Value = "  +0.05%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.914"
$ws.Range("E10").Value = "  -9.26%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.314.58"
$ws.Range("E11").Value = "  -5.90%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "41.30"
$ws.Range("E12").Value = "  -4.45%  "

$ws.Range("E13").Value = "  -4.39%  "

$ws.Range("B14").Value = "WrappedBTC"
$ws.Range("C14").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "93.026.42"
$ws.Range("E14").Value = "  -3.16%  "

$ws.Range("B15").Value = "Toncoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.91"
$ws.Range("E15").Value = "  -4.36%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.931.46"
$ws.Range("E16").Value = "  -5.70%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000241"
$ws.Range("E17").Value = "  -5.35%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.93"
$ws.Range("E18").Value = "  -7.75%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.311.74"
$ws.Range("E19").Value = "  -5.21%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.03"
$ws.Range("E20").Value = "  -8.16%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.84"
$ws.Range("E21").Value = "  -11.61%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.37"
$ws.Range("E22").Value = "  +1.76%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "488.71"
$ws.Range("E23").Value = "  -5.05%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.447"
$ws.Range("E24").Value = "  -10.72%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000179"
$ws.Range("E25").Value = "  -8.63%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.99"
$ws.Range("E26").Value = "  -8.93%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "89.39"
$ws.Range("E27").Value = "  -3.36%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.495.82"
$ws.Range("E28").Value = "  -5.00%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "11.48"
$ws.Range("E29").Value = "  -8.99%  "

$ws.Range("E30").Value = "  +0.24%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "10.92"
$ws.Range("E31").Value = "  -9.32%  "

$ws.Range("E32").Value = "  -4.10%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.62"
$ws.Range("E33").Value = "  -7.04%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.992"
$ws.Range("E34").Value = "  -0.58%  "

$ws.Range("E35").Value = "  -8.62%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "28.11"
$ws.Range("E36").Value = "  -11.23%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.524"
$ws.Range("E37").Value = "  -11.18%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "518.87"
$ws.Range("E38").Value = "  -0.39%  "

$ws.Range("E39").Value = "  -0.10%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.25"
$ws.Range("E40").Value = "  -8.48%  "

$ws.Range("E41").Value = "  -5.00%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.34"
$ws.Range("E42").Value = "  -9.78%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.860"
$ws.Range("E43").Value = "  -9.38%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "24.06"
$ws.Range("E44").Value = "  -0.51%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.65"
$ws.Range("E45").Value = "  -5.56%  "

$ws.Range("B46").Value = "MantraDAO"
$ws.Range("C46").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.55"
$ws.Range("E46").Value = "  -2.66%  "

$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0403"
$ws.Range("E47").Value = "  -4.95%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.29"
$ws.Range("E48").Value = "  -5.66%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "52.57"
$ws.Range("E49").Value = "  -1.68%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.07"
$ws.Range("E50").Value = "  -5.24%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.81"
$ws.Range("E51").Value = "  -6.25%  "
